$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.155.83"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "1.838.19"
$ws.Range("E3").Value = "  -1.41%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6850"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2992"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07420"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07646"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "1.839.43"
$ws.Range("E12").Value = "  -1.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.055"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6802"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "87.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.148"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -7.43%  "
$ws.Range("D17").Value = "29.143.21"
$ws.Range("E17").Value = "  -0.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008153"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.65%  "
$ws.Range("D19").Value = "2.079.38"
$ws.Range("E19").Value = "  -1.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "229.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9996"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.340"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1439"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.698"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.69%  "
$ws.Range("E28").Value = "  -2.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.512"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.261"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.138"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.193"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05259"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7542"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.850"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.133"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.683"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("D38").Value = "1.290.00"
$ws.Range("E38").Value = "  -3.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01824"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.720"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9376"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.930"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "104.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9991"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000124"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.64%  "
$ws.Range("D46").Value = "1.981.44"
$ws.Range("E46").Value = "  -0.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5193"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "64.55"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.479"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.762"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07441"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +17.56%  "
